# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape output (gh-pages generation at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6955
$ws1.Range("F4").Value  = 55
$ws1.Range("F5").Value  = 454
$ws1.Range("F6").Value  = 0
$ws1.Range("F7").Value  = 0
$ws1.Range("F10").Value = 1288
$ws1.Range("F11").Value = 20
$ws1.Range("F12").Value = 108
$ws1.Range("F14").Value = 0
$ws1.Range("F17").Value = 0
$ws1.Range("F20").Value = 5187
$ws1.Range("F22").Value = 0
$ws1.Range("F23").Value = 607
$ws1.Range("F24").Value = 215
$ws1.Range("F25").Value = 224

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 99
$ws4.Range("F7").Value  = 0
$ws4.Range("F8").Value  = 0
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 20
$ws4.Range("F12").Value = 0
$ws4.Range("F17").Value = 48
$ws4.Range("F19").Value = 0
$ws4.Range("F23").Value = 117
$ws4.Range("F25").Value = 0
